$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.247.96'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '3.492.80'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.27'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.06'
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.484'
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.71'
$ws.Range('E9').Value = '  +6.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.387'
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '4.083.38'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '3.488.50'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '64.189.57'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.53'
$ws.Range('E17').Value = '  -5.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.01'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.74'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.47'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.61'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.578'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').Value = '3.630.19'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.39'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.74'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.992'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.24'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.97'
$ws.Range('E32').Value = '  -4.38%  '
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('D34').Value = '3.520.13'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.13'
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.35'
$ws.Range('E37').Value = '  +1.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.85'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.53'
$ws.Range('E39').Value = '  -3.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '163.83'
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0780'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.805'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.37'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.18'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.13'
$ws.Range('E46').Value = '  -6.69%  '
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.929'
$ws.Range('E48').Value = '  +3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.77'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').Value = '2.392.34'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('E51').Value = '  -2.38%  '
